$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sampling-age header row (row 1, columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update "CON" data row (row 2, columns B:E)
$ws.Range("B2").Value = 49.138467275708635
$ws.Range("C2").Value = 55.851088364072623
$ws.Range("D2").Value = 51.940719555127188
$ws.Range("E2").Value = 58.111536687005135

# Update "STR" data row (row 3, columns B:E)
$ws.Range("B3").Value = 44.929400036024902
$ws.Range("C3").Value = 48.857973149811514
$ws.Range("D3").Value = 48.724044589012166
$ws.Range("E3").Value = 56.016889210026186

# Shrink the selected/active range to match the now-trimmed data block
$ws.Range("B1:E3").Select()
